$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.939.64"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").Value = "1.846.59"
$ws.Range("E3").Value = "  +1.23%  "

$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").Value = "'309.39"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").Value = "'0.4771"
$ws.Range("E7").Value = "  +2.74%  "

$ws.Range("D8").Value = "'0.3669"
$ws.Range("E8").Value = "  +1.79%  "

$ws.Range("D9").Value = "'0.07208"
$ws.Range("E9").Value = "  +1.03%  "

$ws.Range("D10").Value = "'0.9275"
$ws.Range("E10").Value = "  +3.03%  "

$ws.Range("E11").Value = "  +1.76%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07691"
$ws.Range("E12").Value = "  -0.91%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.916.35"
$ws.Range("E13").Value = "  +5.04%  "

$ws.Range("D14").Value = "'5.317"
$ws.Range("E14").Value = "  +1.06%  "

$ws.Range("D15").Value = "'6.404"
$ws.Range("E15").Value = "  +1.02%  "

$ws.Range("D16").Value = "'88.87"
$ws.Range("E16").Value = "  +1.20%  "

$ws.Range("E17").Value = "  +0.30%  "

$ws.Range("D18").Value = "'0.000008636"
$ws.Range("E18").Value = "  +0.76%  "

$ws.Range("D19").Value = "'1.009"
$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("D20").Value = "26.975.36"
$ws.Range("E20").Value = "  +1.13%  "

$ws.Range("D21").Value = "'14.55"
$ws.Range("E21").Value = "  +2.71%  "

$ws.Range("D22").Value = "'5.052"
$ws.Range("E22").Value = "  +0.79%  "

$ws.Range("E23").Value = "  +0.92%  "

$ws.Range("D24").Value = "'1.922"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").Value = "'152.59"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("D26").Value = "'18.17"
$ws.Range("E26").Value = "  +1.52%  "

$ws.Range("D27").Value = "'1.996"
$ws.Range("E27").Value = "  +1.36%  "

$ws.Range("D28").Value = "'114.29"
$ws.Range("E28").Value = "  +0.41%  "

$ws.Range("D29").Value = "'4.920"
$ws.Range("E29").Value = "  +1.46%  "

$ws.Range("D30").Value = "'0.08876"
$ws.Range("E30").Value = "  +0.79%  "

$ws.Range("D31").Value = "'3.313"
$ws.Range("E31").Value = "  +5.58%  "

$ws.Range("E32").Value = "  +3.38%  "

$ws.Range("D33").Value = "'0.7450"
$ws.Range("E33").Value = "  +1.75%  "

$ws.Range("D34").Value = "'4.489"
$ws.Range("E34").Value = "  +1.20%  "

$ws.Range("D35").Value = "'2.722"
$ws.Range("E35").Value = "  +0.39%  "

$ws.Range("D36").Value = "'1.107"
$ws.Range("E36").Value = "  +3.00%  "

$ws.Range("D37").Value = "'0.01956"
$ws.Range("E37").Value = "  +1.72%  "

$ws.Range("D38").Value = "'0.05263"
$ws.Range("E38").Value = "  +2.88%  "

$ws.Range("D39").Value = "'2.975"
$ws.Range("E39").Value = "  +1.74%  "

$ws.Range("D40").Value = "'0.5198"
$ws.Range("E40").Value = "  +2.84%  "

$ws.Range("D41").Value = "'6.953"
$ws.Range("E41").Value = "  +0.89%  "

$ws.Range("D42").Value = "'0.1510"
$ws.Range("E42").Value = "  +1.00%  "

$ws.Range("D43").Value = "'8.194"
$ws.Range("E43").Value = "  +2.47%  "

$ws.Range("D44").Value = "'10.51"
$ws.Range("E44").Value = "  +5.52%  "

$ws.Range("D45").Value = "'0.4721"
$ws.Range("E45").Value = "  +1.38%  "

$ws.Range("D46").Value = "'1.010"
$ws.Range("E46").Value = "  +0.25%  "

$ws.Range("D47").Value = "'101.48"
$ws.Range("E47").Value = "  +3.09%  "

$ws.Range("D48").Value = "'1.603"
$ws.Range("E48").Value = "  +2.85%  "

$ws.Range("D49").Value = "'66.03"
$ws.Range("E49").Value = "  +3.48%  "

$ws.Range("D50").Value = "'0.06019"
$ws.Range("E50").Value = "  +0.58%  "

$ws.Range("D51").Value = "'0.8850"
$ws.Range("E51").Value = "  +3.87%  "

Write-Host "Updated cryptos list"
